$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 of data (week 8 entries for PRGE and FALO)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 3

# Update the selected cell to match the post-edit cursor position
$ws.Range("F6").Select() | Out-Null
